$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Cells.Item(2, 10).Value = 4732
$ws.Cells.Item(3, 10).Value = 4998
$ws.Cells.Item(4, 9).Value = 1771
$ws.Cells.Item(4, 10).Value = 1114
$ws.Cells.Item(5, 10).Value = 397
$ws.Cells.Item(6, 10).Value = 6152
$ws.Cells.Item(7, 9).Value = 26219
$ws.Cells.Item(7, 10).Value = 17393

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Cells.Item(2, 10).Value = 131
$ws.Cells.Item(7, 10).Value = 499
$ws.Cells.Item(8, 10).Value = 1120
$ws.Cells.Item(9, 10).Value = 89
$ws.Cells.Item(11, 10).Value = 263
$ws.Cells.Item(18, 10).Value = 152
$ws.Cells.Item(19, 10).Value = 507
$ws.Cells.Item(25, 10).Value = 89
$ws.Cells.Item(29, 10).Value = 993
$ws.Cells.Item(30, 10).Value = 67
$ws.Cells.Item(31, 10).Value = 151
$ws.Cells.Item(33, 10).Value = 795
$ws.Cells.Item(34, 10).Value = 82
$ws.Cells.Item(36, 10).Value = 242
$ws.Cells.Item(37, 10).Value = 549
$ws.Cells.Item(40, 10).Value = 40
$ws.Cells.Item(42, 10).Value = 703
$ws.Cells.Item(43, 10).Value = 147
$ws.Cells.Item(48, 10).Value = 189
$ws.Cells.Item(51, 10).Value = 219
$ws.Cells.Item(52, 10).Value = 440
$ws.Cells.Item(53, 10).Value = 218
$ws.Cells.Item(54, 10).Value = 340
$ws.Cells.Item(57, 10).Value = 74
$ws.Cells.Item(60, 10).Value = 113
$ws.Cells.Item(63, 9).Value = 230
$ws.Cells.Item(63, 10).Value = 67
$ws.Cells.Item(64, 10).Value = 118
$ws.Cells.Item(65, 10).Value = 464
$ws.Cells.Item(66, 10).Value = 55
$ws.Cells.Item(67, 10).Value = 670
$ws.Cells.Item(68, 10).Value = 31
$ws.Cells.Item(72, 10).Value = 68
$ws.Cells.Item(75, 10).Value = 51
$ws.Cells.Item(76, 9).Value = 381
$ws.Cells.Item(76, 10).Value = 250
$ws.Cells.Item(77, 10).Value = 131
$ws.Cells.Item(78, 10).Value = 218
$ws.Cells.Item(79, 10).Value = 500
$ws.Cells.Item(83, 10).Value = 378
$ws.Cells.Item(85, 10).Value = 770
$ws.Cells.Item(86, 9).Value = 169
$ws.Cells.Item(87, 10).Value = 61
$ws.Cells.Item(89, 10).Value = 222
$ws.Cells.Item(91, 10).Value = 193
$ws.Cells.Item(95, 10).Value = 262
$ws.Cells.Item(96, 10).Value = 210
$ws.Cells.Item(98, 10).Value = 110
$ws.Cells.Item(99, 10).Value = 259
$ws.Cells.Item(101, 9).Value = 26219
$ws.Cells.Item(101, 10).Value = 17393

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Cells.Item(6, 10).Value = 158
$ws.Cells.Item(7, 10).Value = 499

$ws = $wb.Worksheets.Item("Uptown")
$ws.Cells.Item(6, 10).Value = 66
$ws.Cells.Item(7, 10).Value = 222

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Cells.Item(2, 10).Value = 89
$ws.Cells.Item(7, 10).Value = 263

$ws = $wb.Worksheets.Item("Little Village")
$ws.Cells.Item(2, 10).Value = 104
$ws.Cells.Item(6, 10).Value = 189
$ws.Cells.Item(7, 10).Value = 440

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Cells.Item(5, 10).Value = 4
$ws.Cells.Item(6, 10).Value = 72
$ws.Cells.Item(7, 10).Value = 210

$ws = $wb.Worksheets.Item("Austin")
$ws.Cells.Item(4, 10).Value = 59
$ws.Cells.Item(6, 10).Value = 372
$ws.Cells.Item(7, 10).Value = 1120

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Cells.Item(2, 10).Value = 46
$ws.Cells.Item(7, 10).Value = 218

$ws = $wb.Worksheets.Item("South Shore")
$ws.Cells.Item(3, 10).Value = 279
$ws.Cells.Item(6, 10).Value = 219
$ws.Cells.Item(7, 10).Value = 770

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Cells.Item(2, 10).Value = 114
$ws.Cells.Item(7, 10).Value = 378

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Cells.Item(6, 10).Value = 64
$ws.Cells.Item(7, 10).Value = 262

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Cells.Item(2, 10).Value = 72
$ws.Cells.Item(6, 10).Value = 71
$ws.Cells.Item(7, 10).Value = 259

$ws = $wb.Worksheets.Item("Fuller Park")
$ws.Cells.Item(6, 10).Value = 18
$ws.Cells.Item(7, 10).Value = 67

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Cells.Item(2, 10).Value = 164
$ws.Cells.Item(7, 10).Value = 670

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Cells.Item(6, 10).Value = 39
$ws.Cells.Item(7, 10).Value = 151

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Cells.Item(2, 10).Value = 164
$ws.Cells.Item(7, 10).Value = 549

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Cells.Item(2, 10).Value = 198
$ws.Cells.Item(6, 10).Value = 275
$ws.Cells.Item(7, 10).Value = 795

$ws = $wb.Worksheets.Item("New City")
$ws.Cells.Item(3, 10).Value = 140
$ws.Cells.Item(5, 10).Value = 12
$ws.Cells.Item(7, 10).Value = 464

$ws = $wb.Worksheets.Item("Loop")
$ws.Cells.Item(6, 10).Value = 160
$ws.Cells.Item(7, 10).Value = 340

$ws = $wb.Worksheets.Item("Englewood")
$ws.Cells.Item(3, 10).Value = 343
$ws.Cells.Item(6, 10).Value = 255
$ws.Cells.Item(7, 10).Value = 993

$ws = $wb.Worksheets.Item("Chatham")
$ws.Cells.Item(3, 10).Value = 145
$ws.Cells.Item(6, 10).Value = 188
$ws.Cells.Item(7, 10).Value = 507

$ws = $wb.Worksheets.Item("Lake View")
$ws.Cells.Item(6, 10).Value = 95
$ws.Cells.Item(7, 10).Value = 189

$ws = $wb.Worksheets.Item("River North")
$ws.Cells.Item(4, 10).Value = 23
$ws.Cells.Item(5, 9).Value = 8
$ws.Cells.Item(7, 9).Value = 381
$ws.Cells.Item(7, 10).Value = 250

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Cells.Item(2, 10).Value = 150
$ws.Cells.Item(6, 10).Value = 360
$ws.Cells.Item(7, 10).Value = 703

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Cells.Item(3, 10).Value = 76
$ws.Cells.Item(7, 10).Value = 218

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Cells.Item(2, 10).Value = 60
$ws.Cells.Item(6, 10).Value = 37
$ws.Cells.Item(7, 10).Value = 193

$ws = $wb.Worksheets.Item("Roseland")
$ws.Cells.Item(2, 10).Value = 142
$ws.Cells.Item(3, 10).Value = 176
$ws.Cells.Item(6, 10).Value = 137
$ws.Cells.Item(7, 10).Value = 500

$ws = $wb.Worksheets.Item("Near South Side")
$ws.Cells.Item(6, 10).Value = 41
$ws.Cells.Item(7, 10).Value = 118

$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Cells.Item(6, 10).Value = 78
$ws.Cells.Item(7, 10).Value = 152

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Cells.Item(6, 10).Value = 70
$ws.Cells.Item(7, 10).Value = 242

$ws = $wb.Worksheets.Item("Garfield Ridge")
$ws.Cells.Item(2, 10).Value = 25
$ws.Cells.Item(7, 10).Value = 82

$ws = $wb.Worksheets.Item("East Side")
$ws.Cells.Item(3, 10).Value = 28
$ws.Cells.Item(7, 10).Value = 89

$ws = $wb.Worksheets.Item("Wicker Park")
$ws.Cells.Item(6, 10).Value = 63
$ws.Cells.Item(7, 10).Value = 110

$ws = $wb.Worksheets.Item("North Center")
$ws.Cells.Item(6, 10).Value = 31
$ws.Cells.Item(7, 10).Value = 55

$ws = $wb.Worksheets.Item("Avalon Park")
$ws.Cells.Item(2, 10).Value = 25
$ws.Cells.Item(7, 10).Value = 89

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Cells.Item(6, 10).Value = 47
$ws.Cells.Item(7, 10).Value = 131

$ws = $wb.Worksheets.Item("Streeterville")
$ws.Cells.Item(5, 9).Value = 3
$ws.Cells.Item(7, 9).Value = 169

$ws = $wb.Worksheets.Item("Pullman")
$ws.Cells.Item(2, 10).Value = 24
$ws.Cells.Item(7, 10).Value = 51

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Cells.Item(6, 10).Value = 77
$ws.Cells.Item(7, 10).Value = 219

$ws = $wb.Worksheets.Item("North Park")
$ws.Cells.Item(2, 10).Value = 14
$ws.Cells.Item(7, 10).Value = 31

$ws = $wb.Worksheets.Item("Mckinley Park")
$ws.Cells.Item(6, 10).Value = 27
$ws.Cells.Item(7, 10).Value = 74

$ws = $wb.Worksheets.Item("Morgan Park")
$ws.Cells.Item(3, 10).Value = 33
$ws.Cells.Item(7, 10).Value = 113

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Cells.Item(3, 10).Value = 28
$ws.Cells.Item(6, 10).Value = 86
$ws.Cells.Item(7, 10).Value = 147

$ws = $wb.Worksheets.Item("Old Town")
$ws.Cells.Item(6, 10).Value = 22
$ws.Cells.Item(7, 10).Value = 68

$ws = $wb.Worksheets.Item("Riverdale")
$ws.Cells.Item(3, 10).Value = 45
$ws.Cells.Item(7, 10).Value = 131

$ws = $wb.Worksheets.Item("Hegewisch")
$ws.Cells.Item(4, 10).Value = 4
$ws.Cells.Item(7, 10).Value = 40

$ws = $wb.Worksheets.Item("Ukrainian Village")
$ws.Cells.Item(6, 10).Value = 38
$ws.Cells.Item(7, 10).Value = 61
